$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change Mess Kit quantity from 1 to 0
$ws.Range("B4").Value = 0

# Add new row 7: Portable Radio
$ws.Range("A7").Value = "Portable Radio"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Comes with AM/FM and a pretty good volume for something so small"

# Copy formatting from row 6 (existing similarly-styled row) into the new row 7 cells
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null

$ws.Range("D6").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update active selection from E7 to E8
$ws.Range("E8").Select() | Out-Null
